$wb = $excel.ActiveWorkbook

# Rows (by sheet row number) that were re-handed-off and need their
# "Priority" marked as "ht" plus a refreshed handoff timestamp.
$rows = @(7, 10, 11, 12, 13, 14)

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsDeDe = $wb.Worksheets.Item("de-de")

foreach ($r in $rows) {
    $wsOverview.Range("G$r").Value = "2016-08-18 14:22:14"

    $wsZhCn.Range("E$r").Value = "ht"
    $wsZhCn.Range("H$r").Value = "2016-08-18 14:22:00"

    $wsDeDe.Range("E$r").Value = "ht"
    $wsDeDe.Range("H$r").Value = "2016-08-18 14:22:14"
}
